$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

Set-TextCell "D2" "301.61"
Set-TextCell "E2" "-2.67%"
Set-TextCell "D3" "35.33"
Set-TextCell "E3" "-0.54%"
Set-TextCell "D4" "5.074"
Set-TextCell "E4" "-0.24%"
Set-TextCell "D5" "0.07925"
Set-TextCell "E5" "-2.90%"
Set-TextCell "D6" "1.889"
Set-TextCell "E6" "-8.02%"
Set-TextCell "B7" "KuCoinToken"
Set-TextCell "C7" "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextCell "D7" "7.781"
Set-TextCell "E7" "-2.01%"
Set-TextCell "B8" "GateToken"
Set-TextCell "C8" "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextCell "D8" "4.047"
Set-TextCell "E8" "-2.02%"
Set-TextCell "D9" "0.9293"
Set-TextCell "D10" "0.1382"
Set-TextCell "E10" "30.00%"
Set-TextCell "D11" "0.1899"
Set-TextCell "E11" "-1.12%"
Set-TextCell "D12" "0.09135"
Set-TextCell "E12" "-1.79%"
Set-TextCell "D13" "0.03437"
Set-TextCell "E13" "-6.51%"
Set-TextCell "D14" "0.09834"
Set-TextCell "E14" "-0.52%"
Set-TextCell "D15" "0.001396"
Set-TextCell "E15" "-2.46%"
Set-TextCell "D16" "0.005881"
Set-TextCell "E16" "3.39%"
Set-TextCell "D17" "3.532"
Set-TextCell "E17" "1.59%"
Set-TextCell "E18" "1.44%"
Set-TextCell "D19" "0.3424"
Set-TextCell "E19" "0.95%"
Set-TextCell "D20" "0.1304"
Set-TextCell "E20" "0.31%"
Set-TextCell "E21" "-1.01%"
Set-TextCell "D22" "0.2402"
Set-TextCell "E22" "8.57%"
Set-TextCell "D23" "0.04498"
Set-TextCell "E23" "-1.18%"
Set-TextCell "D24" "0.001214"
Set-TextCell "E24" "-1.11%"
Set-TextCell "D25" "0.004759"
Set-TextCell "E25" "-0.45%"
Set-TextCell "D26" "0.0001231"
Set-TextCell "E26" "-1.63%"
Set-TextCell "D27" "0.0003003"
Set-TextCell "E27" "-32.50%"
Set-TextCell "D39" "0.01857"
Set-TextCell "E39" "-5.40%"
Set-TextCell "D40" "0.04764"
Set-TextCell "E40" "-2.42%"
Set-TextCell "D41" "0.007364"
Set-TextCell "E41" "-2.54%"
Set-TextCell "D42" "0.009613"
Set-TextCell "E42" "-2.90%"
Set-TextCell "E43" "-4.09%"
Set-TextCell "D44" "0.002112"
Set-TextCell "E44" "-4.98%"
Set-TextCell "D45" "0.01088"
Set-TextCell "E45" "-6.38%"
Set-TextCell "D46" "0.00006244"
Set-TextCell "E46" "-5.63%"
Set-TextCell "E47" "0.07%"
Set-TextCell "D48" "64.68"
Set-TextCell "E48" "7.75%"
Set-TextCell "E49" "10.59%"
Set-TextCell "D50" "0.00002102"
Set-TextCell "E50" "0.07%"
Set-TextCell "D51" "0.0002002"
Set-TextCell "E51" "0.07%"
